$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 2.543891666666667
$ws.Cells.Item(2, 8).Value = 7.631675
$ws.Cells.Item(2, 9).Value = 0.0487891203504456
$ws.Cells.Item(2, 10).Value = 0.04878912035044559
$ws.Cells.Item(2, 13).Value = 0.09834766666666667
$ws.Cells.Item(2, 14).Value = 0.295043
$ws.Cells.Item(2, 15).Value = 0.2818566198948398
$ws.Cells.Item(2, 16).Value = 0.2818566198948398
$ws.Cells.Item(2, 17).Value = 0.2501858096694444
$ws.Cells.Item(2, 18).Value = 2.251672287025
$ws.Cells.Item(2, 19).Value = 0.01375153654961914
$ws.Cells.Item(2, 20).Value = 0.01375153654961914

# Row 3
$ws.Cells.Item(3, 7).Value = 2.543891666666667
$ws.Cells.Item(3, 8).Value = 7.631675
$ws.Cells.Item(3, 9).Value = 0.0487891203504456
$ws.Cells.Item(3, 10).Value = 0.04878912035044559
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.2505803333333333
$ws.Cells.Item(3, 14).Value = 0.751741
$ws.Cells.Item(3, 15).Value = 0.7181433801051602
$ws.Cells.Item(3, 16).Value = 0.7181433801051602
$ws.Cells.Item(3, 17).Value = 0.6374492217972223
$ws.Cells.Item(3, 18).Value = 5.737042996175
$ws.Cells.Item(3, 19).Value = 0.03503758380082646
$ws.Cells.Item(3, 20).Value = 0.03503758380082646

# Row 4
$ws.Cells.Item(4, 9).Value = 0.0400662233111763
$ws.Cells.Item(4, 10).Value = 0.0400662233111763
$ws.Cells.Item(4, 13).Value = 0.09834766666666667
$ws.Cells.Item(4, 14).Value = 0.295043
$ws.Cells.Item(4, 15).Value = 0.2818566198948398
$ws.Cells.Item(4, 16).Value = 0.2818566198948398
$ws.Cells.Item(4, 17).Value = 0.2054556517416666
$ws.Cells.Item(4, 18).Value = 1.849100865675
$ws.Cells.Item(4, 19).Value = 0.01129293027443999
$ws.Cells.Item(4, 20).Value = 0.01129293027443999

# Row 5
$ws.Cells.Item(5, 9).Value = 0.0400662233111763
$ws.Cells.Item(5, 10).Value = 0.0400662233111763
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.2505803333333333
$ws.Cells.Item(5, 14).Value = 0.751741
$ws.Cells.Item(5, 15).Value = 0.7181433801051602
$ws.Cells.Item(5, 16).Value = 0.7181433801051602
$ws.Cells.Item(5, 17).Value = 0.5234811098583333
$ws.Cells.Item(5, 18).Value = 4.711329988725
$ws.Cells.Item(5, 19).Value = 0.02877329303673631
$ws.Cells.Item(5, 20).Value = 0.02877329303673632

# Row 6
$ws.Cells.Item(6, 7).Value = 1.800112666666666
$ws.Cells.Item(6, 8).Value = 5.400338
$ws.Cells.Item(6, 9).Value = 0.03452423493074386
$ws.Cells.Item(6, 10).Value = 0.03452423493074386
$ws.Cells.Item(6, 13).Value = 0.09834766666666667
$ws.Cells.Item(6, 14).Value = 0.295043
$ws.Cells.Item(6, 15).Value = 0.2818566198948398
$ws.Cells.Item(6, 16).Value = 0.2818566198948398
$ws.Cells.Item(6, 17).Value = 0.1770368805037777
$ws.Cells.Item(6, 18).Value = 1.593331924534
$ws.Cells.Item(6, 19).Value = 0.009730884162034823
$ws.Cells.Item(6, 20).Value = 0.009730884162034825

# Row 7
$ws.Cells.Item(7, 7).Value = 1.800112666666666
$ws.Cells.Item(7, 8).Value = 5.400338
$ws.Cells.Item(7, 9).Value = 0.03452423493074386
$ws.Cells.Item(7, 10).Value = 0.03452423493074386
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.2505803333333333
$ws.Cells.Item(7, 14).Value = 0.751741
$ws.Cells.Item(7, 15).Value = 0.7181433801051602
$ws.Cells.Item(7, 16).Value = 0.7181433801051602
$ws.Cells.Item(7, 17).Value = 0.4510728320508889
$ws.Cells.Item(7, 18).Value = 4.059655488458
$ws.Cells.Item(7, 19).Value = 0.02479335076870904
$ws.Cells.Item(7, 20).Value = 0.02479335076870904

# Row 8
$ws.Cells.Item(8, 7).Value = 45.70747266666666
$ws.Cells.Item(8, 8).Value = 137.122418
$ws.Cells.Item(8, 9).Value = 0.8766204214076343
$ws.Cells.Item(8, 10).Value = 0.8766204214076342
$ws.Cells.Item(8, 13).Value = 0.09834766666666667
$ws.Cells.Item(8, 14).Value = 0.295043
$ws.Cells.Item(8, 15).Value = 0.2818566198948398
$ws.Cells.Item(8, 16).Value = 0.2818566198948398
$ws.Cells.Item(8, 17).Value = 4.495223285997111
$ws.Cells.Item(8, 18).Value = 40.45700957397399
$ws.Cells.Item(8, 19).Value = 0.2470812689087458
$ws.Cells.Item(8, 20).Value = 0.2470812689087459

# Row 9
$ws.Cells.Item(9, 7).Value = 45.70747266666666
$ws.Cells.Item(9, 8).Value = 137.122418
$ws.Cells.Item(9, 9).Value = 0.8766204214076343
$ws.Cells.Item(9, 10).Value = 0.8766204214076342
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.2505803333333333
$ws.Cells.Item(9, 14).Value = 0.751741
$ws.Cells.Item(9, 15).Value = 0.7181433801051602
$ws.Cells.Item(9, 16).Value = 0.7181433801051602
$ws.Cells.Item(9, 17).Value = 11.45339373663755
$ws.Cells.Item(9, 18).Value = 103.080543629738
$ws.Cells.Item(9, 19).Value = 0.6295391524988885
$ws.Cells.Item(9, 20).Value = 0.6295391524988884
